$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1598.75
$ws.Range("I12").Value = 1731.6666
$ws.Range("K12").Value = 1731.6666
$ws.Range("M12").Value = -1561.6666
$ws.Range("H17").Value = 640173.1
$ws.Range("J17").Value = 640173.1
$ws.Range("L17").Value = 1920519.3
$ws.Range("N17").Value = -1920855.3
$ws.Range("H64").Value = 11849.111
$ws.Range("I64").Value = 3632.6667
$ws.Range("J64").Value = 15957.333
$ws.Range("K64").Value = 3632.6667
$ws.Range("L64").Value = 15957.333
$ws.Range("M64").Value = -3384.6667
$ws.Range("N64").Value = -16453.333
$ws.Range("H67").Value = 11849.111
$ws.Range("I67").Value = 3632.6667
$ws.Range("J67").Value = 15957.333
$ws.Range("K67").Value = 3632.6667
$ws.Range("L67").Value = 15957.333
$ws.Range("M67").Value = -2774.6667
$ws.Range("N67").Value = -17673.333
$ws.Range("H113").Value = 64544.473
$ws.Range("I113").Value = 135749.88
$ws.Range("J113").Value = 12758.728
$ws.Range("K113").Value = 135749.88
$ws.Range("L113").Value = 12758.728
$ws.Range("M113").Value = -132495.88
$ws.Range("N113").Value = -19266.728
$ws.Range("H116").Value = 11998
$ws.Range("J116").Value = 9166.5
$ws.Range("L116").Value = 9166.5
$ws.Range("N116").Value = -16050.5
$ws.Range("H132").Value = 2321.375
$ws.Range("I132").Value = 1403.8182
$ws.Range("J132").Value = 4340
$ws.Range("K132").Value = 4211.4546
$ws.Range("L132").Value = 13020
$ws.Range("M132").Value = -1681.4546
$ws.Range("N132").Value = -18080
$ws.Range("H138").Value = 3455.3425
$ws.Range("I138").Value = 1848.7307
$ws.Range("K138").Value = 5546.1921
$ws.Range("M138").Value = -406.1921000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 255499.75
$ws.Range("I2").Value = 338166.34
$ws.Range("K2").Value = 338166.34
$ws.Range("M2").Value = -338053.34
$ws.Range("H6").Value = 217.5
$ws.Range("I6").Value = 217.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 217.5
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -44.5
$ws.Range("H45").Value = 140854.12
$ws.Range("I45").Value = 185483.17
$ws.Range("K45").Value = 185483.17
$ws.Range("M45").Value = -185106.17
$ws.Range("H74").Value = 13033.25
$ws.Range("I74").Value = 1903.4546
$ws.Range("J74").Value = 37518.8
$ws.Range("K74").Value = 1903.4546
$ws.Range("L74").Value = 37518.8
$ws.Range("M74").Value = -1029.4546
$ws.Range("N74").Value = -39266.8
$ws.Range("H77").Value = 13033.25
$ws.Range("I77").Value = 1903.4546
$ws.Range("J77").Value = 37518.8
$ws.Range("K77").Value = 9517.273000000001
$ws.Range("L77").Value = 187594
$ws.Range("M77").Value = -5149.273000000001
$ws.Range("N77").Value = -196330
$ws.Range("H97").Value = 536.9375
$ws.Range("I97").Value = 499.35715
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 499.35715
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -3.35714999999999
$ws.Range("N97").Value = -1792
$ws.Range("H102").Value = 4835.2856
$ws.Range("I102").Value = 3769.6
$ws.Range("J102").Value = 7499.5
$ws.Range("K102").Value = 3769.6
$ws.Range("L102").Value = 7499.5
$ws.Range("M102").Value = -2147.6
$ws.Range("N102").Value = -10743.5
$ws.Range("H116").Value = 255499.75
$ws.Range("I116").Value = 338166.34
$ws.Range("K116").Value = 338166.34
$ws.Range("M116").Value = -335872.34
$ws.Range("H132").Value = 3709.8845
$ws.Range("I132").Value = 3252.2222
$ws.Range("K132").Value = 9756.6666
$ws.Range("M132").Value = -7226.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 255499.75
$ws.Range("I3").Value = 338166.34
$ws.Range("K3").Value = 338166.34
$ws.Range("M3").Value = -338052.34
$ws.Range("H86").Value = 2076.923
$ws.Range("H89").Value = 2076.923
$ws.Range("H94").Value = 1410
$ws.Range("I94").Value = 1395
$ws.Range("K94").Value = 1395
$ws.Range("M94").Value = -944
$ws.Range("H134").Value = 1735.9259
$ws.Range("I134").Value = 1725.7693
$ws.Range("K134").Value = 5177.3079
$ws.Range("M134").Value = -2642.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 642.8570999999999
$ws.Range("I22").Value = 640
$ws.Range("K22").Value = 640
$ws.Range("M22").Value = -290
$ws.Range("H31").Value = 19128.508
$ws.Range("I31").Value = 27230.05
$ws.Range("K31").Value = 27230.05
$ws.Range("M31").Value = -26935.05
$ws.Range("H34").Value = 19128.508
$ws.Range("I34").Value = 27230.05
$ws.Range("K34").Value = 27230.05
$ws.Range("M34").Value = -27028.05
$ws.Range("H58").Value = 2379.7144
$ws.Range("I58").Value = 2475.25
$ws.Range("K58").Value = 2475.25
$ws.Range("M58").Value = -2272.25
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").ClearContents()
$ws.Range("N115").Value = 0
$ws.Range("H121").Value = 14296
$ws.Range("I121").Value = 14296
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 14296
$ws.Range("L121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -12986
$ws.Range("H132").Value = 5408.604
$ws.Range("I132").Value = 3707.756
$ws.Range("K132").Value = 11123.268
$ws.Range("M132").Value = -8593.268
$ws.Range("H136").Value = 2379.7144
$ws.Range("I136").Value = 2475.25
$ws.Range("K136").Value = 7425.75
$ws.Range("M136").Value = -4875.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1113.7059
$ws.Range("I132").Value = 1076.9678
$ws.Range("J132").Value = 1493.3334
$ws.Range("K132").Value = 9692.7102
$ws.Range("L132").Value = 13440.0006
$ws.Range("M132").Value = -7162.7102
$ws.Range("N132").Value = -18500.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 22500
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H48").Value = 19000
$ws.Range("I48").Value = 19000
$ws.Range("K48").Value = 19000
$ws.Range("M48").Value = -18515
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = 0
$ws.Range("H97").Value = 1037
$ws.Range("I97").Value = 1588.875
$ws.Range("K97").Value = 1588.875
$ws.Range("M97").Value = -1092.875
$ws.Range("H102").Value = 100002910
$ws.Range("I102").Value = 2683.1667
$ws.Range("J102").Value = 250003250
$ws.Range("K102").Value = 2683.1667
$ws.Range("L102").Value = 250003250
$ws.Range("M102").Value = -1061.1667
$ws.Range("N102").Value = -250006494
$ws.Range("H132").Value = 5201.6665
$ws.Range("I132").Value = 4688
$ws.Range("J132").Value = 6999.5
$ws.Range("K132").Value = 14064
$ws.Range("L132").Value = 20998.5
$ws.Range("M132").Value = -11534
$ws.Range("N132").Value = -26058.5
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 11746.25
$ws.Range("I43").Value = 8995
$ws.Range("K43").Value = 8995
$ws.Range("M43").Value = -8802
$ws.Range("H68").Value = 3010.7693
$ws.Range("I68").Value = 2420
$ws.Range("K68").Value = 2420
$ws.Range("M68").Value = -1671
$ws.Range("H71").Value = 3010.7693
$ws.Range("I71").Value = 2420
$ws.Range("K71").Value = 12100
$ws.Range("M71").Value = -8356
$ws.Range("H132").Value = 6499.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2240.7693
$ws.Range("I96").Value = 2372.4
$ws.Range("K96").Value = 2372.4
$ws.Range("M96").Value = -999.4000000000001
$ws.Range("H132").Value = 1839.88
$ws.Range("I132").Value = 1885.5714
$ws.Range("K132").Value = 5656.7142
$ws.Range("M132").Value = -3126.7142
$ws.Range("H136").Value = 1871.804
$ws.Range("I136").Value = 1510.5
$ws.Range("J136").Value = 2387.9524
$ws.Range("K136").Value = 4531.5
$ws.Range("L136").Value = 7163.8572
$ws.Range("M136").Value = -1981.5
$ws.Range("N136").Value = -12263.8572
